$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = -6
